$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename "Vendas" -> "Marco", drop "Mes Venda" column (G),
#     and keep only the rows that belong to March (original rows 2 and 4).
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Março"
$ws1.Columns.Item(7).Delete()   # drop "Mes Venda" column
$ws1.Rows.Item(5).Delete()      # original row5 -> Abril (moved out)
$ws1.Rows.Item(3).Delete()      # original row3 -> Janeiro (moved out)

# --- Sheet 2: "Janeiro" (former row 3: devic / Cimento nassau / 23 / 1 / pix / 23)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Janeiro"
$ws2.Range("A1").Value = "Cliente"
$ws2.Range("B1").Value = "Produto"
$ws2.Range("C1").Value = "Valor da venda"
$ws2.Range("D1").Value = "Quantidade"
$ws2.Range("E1").Value = "Pagamentos"
$ws2.Range("F1").Value = "Dia Venda"
$ws2.Range("A2").Value = "devic"
$ws2.Range("B2").Value = "Cimento nassau"
$ws2.Range("C2").Value = 23.0
$ws2.Range("D2").Value = 1
$ws2.Range("E2").Value = "pix"
$ws2.Range("F2").Value = 23

# --- Sheet 3: "Abril" (former row 5: devic / Cimento nassau / 35 / 1 / pix / 29)
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Abril"
$ws3.Range("A1").Value = "Cliente"
$ws3.Range("B1").Value = "Produto"
$ws3.Range("C1").Value = "Valor da venda"
$ws3.Range("D1").Value = "Quantidade"
$ws3.Range("E1").Value = "Pagamentos"
$ws3.Range("F1").Value = "Dia Venda"
$ws3.Range("A2").Value = "devic"
$ws3.Range("B2").Value = "Cimento nassau"
$ws3.Range("C2").Value = 35.0
$ws3.Range("D2").Value = 1
$ws3.Range("E2").Value = "pix"
$ws3.Range("F2").Value = 29

$ws1.Select()
